$wb = $excel.ActiveWorkbook

# Hitos sheet
$ws1 = $wb.Worksheets.Item("Hitos")
$ws1.Columns.Item(2).ColumnWidth = 23.666666666666668
$ws1.Range("A7:C9").Select() | Out-Null
$ws1.Application.ActiveWindow.Zoom = 115

# Tareas divididas sheet
$ws2 = $wb.Worksheets.Item("Tareas divididas")

# Clear old data row 1
$ws2.Range("A1:B1").ClearContents()

# Shift data: row2 = Leandro, row3 = Ezequiel, row4 = Franco
$ws2.Range("A2").Value = "Leandro"
$ws2.Range("A3").Value = "Ezequiel"
$ws2.Range("B3").ClearContents()
$ws2.Range("A4").Value = "Franco"
$ws2.Range("B4").Value = "Stock, Bajas logicas de categorias y menus"
$ws2.Range("B2").Value = "Reportes, Armado procedimiento de suma del total a pagar, abm mesas"

$ws2.Columns.Item(2).ColumnWidth = 64.5

$ws2.Activate() | Out-Null
$ws2.Range("A7:B10").Select() | Out-Null
